$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete last data row (row 53); the series shifted up by
# one period and the oldest observation was dropped.
$ws.Rows("53:53").Delete()

# Build the full corrected data block (rows 2-52, columns A-E) reflecting the
# bugfixed naive forecaster output (dates shifted, recomputed C/E values).
$data = @(
    @(39583,2008,"",2009,""),
    @(39765,2008,"",2009,""),
    @(39948,2009,"",2010,""),
    @(40130,2009,1.834695583582535,2010,""),
    @(40310,2010,"",2011,""),
    @(40494,2010,1.767835936772144,2011,""),
    @(40676,2011,"",2012,""),
    @(40862,2011,1.074400434091038,2012,""),
    @(41044,2012,"",2013,""),
    @(41228,2012,0.9212998022035679,2013,1.274704633957136),
    @(41409,2013,1.097054137926201,2014,1.37755776875883),
    @(41592,2013,1.141837882844188,2014,1.404348988410131),
    @(41774,2014,1.385527545913412,2015,1.329814931661888),
    @(41957,2014,1.335361538769475,2015,1.269653854937691),
    @(42137,2015,1.296301936385214,2016,1.355477993452414),
    @(42321,2015,1.202048372526998,2016,1.253742200752095),
    @(42503,2016,2.441628883342295,2017,1.386547975635688),
    @(42689,2016,2.677488680362305,2017,1.805615391969595),
    @(42867,2017,2.565764046666463,2018,1.833587970352424),
    @(43053,2017,2.466954516646402,2018,1.661541796722577),
    @(43145,2018,1.376993627314671,2019,1.553425185274571),
    @(43235,2018,1.263447557103259,2019,1.485511920344451),
    @(43326,2018,1.362852986880547,2019,1.635045928803081),
    @(43418,2018,1.401189216021326,2019,1.815016201748643),
    @(43510,2019,1.868356483387124,2020,1.600603483732033),
    @(43600,2019,2.117022522597423,2020,1.745834498329324),
    @(43691,2019,2.244229492550187,2020,1.944926006147352),
    @(43783,2019,2.217567799050979,2020,1.810449264563152),
    @(43875,2020,2.005372766276703,2021,1.679039075934385),
    @(43966,2020,2.149400276001101,2021,1.76475225558832),
    @(44068,2020,2.139672475020404,2021,1.674992401025355),
    @(44159,2020,2.139672475020404,2021,2.128328071999674),
    @(44251,2021,2.273520076663971,2022,1.827015901454287),
    @(44341,2021,2.453568910971748,2022,2.131436976903012),
    @(44432,2021,2.100991693542231,2022,0.5797569954799853),
    @(44525,2021,2.100991693542231,2022,1.11435041103376),
    @(44617,2022,0.7380952286421882,2023,1.421947874507667),
    @(44706,2022,0.812682184439506,2023,1.556352278772266),
    @(44798,2022,0.8967077601845341,2023,1.891749670939347),
    @(44890,2022,0.8967077601845341,2023,0.3338851812143995),
    @(44981,2023,-0.4131792716363547,2024,1.173782198617435),
    @(45071,2023,0.9940067218177528,2024,1.820779918499094),
    @(45163,2023,0.782207885866093,2024,1.655852464312013),
    @(45254,2023,0.782207885866093,2024,2.228542839642689),
    @(45345,2024,1.805571054927801,2025,1.632015075917925),
    @(45436,2024,1.634555928116921,2025,1.554016159863814),
    @(45534,2024,1.508385007449875,2025,1.310740881193517),
    @(45618,2024,1.508385007449875,2025,0.9823016603409229),
    @(45713,2025,0.6395223689078522,2026,1.590970324046337),
    @(45800,2025,0.6231570351797,2026,1.581524829939718),
    @(45891,2025,0.6014263374495288,2026,1.539180932382078)
)

$numRows = $data.Length
$arr = New-Object 'object[,]' $numRows,5
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt 5; $j++) {
        $val = $data[$i][$j]
        if ($val -eq "") {
            $arr[$i,$j] = $null
        } else {
            $arr[$i,$j] = $val
        }
    }
}

$rng = $ws.Range("A2:E52")
$rng.Value2 = $arr
